$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "週末特価" post row (row 66) entirely; Excel shifts all
# subsequent rows (67-180) up by one, producing rows 66-179.
$ws.Rows("66").Delete()
